$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 38750
$ws.Range("I18").Value = 38333.332
$ws.Range("J18").Value = 40000
$ws.Range("K18").Value = 38333.332
$ws.Range("L18").Value = 40000
$ws.Range("M18").Value = -38049.332
$ws.Range("N18").Value = -40568
$ws.Range("H31").Value = 341.33334
$ws.Range("I31").Value = 341.33334
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1024.00002
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -794.0000199999999
$ws.Range("H40").Value = 3644.4443
$ws.Range("I40").Value = 2660
$ws.Range("J40").Value = 4875
$ws.Range("K40").Value = 2660
$ws.Range("L40").Value = 4875
$ws.Range("M40").Value = -2485
$ws.Range("N40").Value = -5225
$ws.Range("H41").Value = 1223.3
$ws.Range("I41").Value = 1091
$ws.Range("J41").Value = 1532
$ws.Range("K41").Value = 1091
$ws.Range("L41").Value = 1532
$ws.Range("M41").Value = -651
$ws.Range("N41").Value = -2412
$ws.Range("H46").Value = 1668770.4
$ws.Range("I46").Value = 1800
$ws.Range("J46").Value = 2002164.4
$ws.Range("K46").Value = 5400
$ws.Range("L46").Value = 6006493.199999999
$ws.Range("M46").Value = -5281
$ws.Range("N46").Value = -6006731.199999999
$ws.Range("H60").Value = 1668770.4
$ws.Range("I60").Value = 1800
$ws.Range("J60").Value = 2002164.4
$ws.Range("K60").Value = 5400
$ws.Range("L60").Value = 6006493.199999999
$ws.Range("M60").Value = -4916
$ws.Range("N60").Value = -6007461.199999999
$ws.Range("H64").Value = 4117.9546
$ws.Range("I64").Value = 3733.0557
$ws.Range("J64").Value = 5850
$ws.Range("K64").Value = 3733.0557
$ws.Range("L64").Value = 5850
$ws.Range("M64").Value = -3485.0557
$ws.Range("N64").Value = -6346
$ws.Range("H67").Value = 4117.9546
$ws.Range("I67").Value = 3733.0557
$ws.Range("J67").Value = 5850
$ws.Range("K67").Value = 3733.0557
$ws.Range("L67").Value = 5850
$ws.Range("M67").Value = -2875.0557
$ws.Range("N67").Value = -7566
$ws.Range("H68").Value = 25295
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 25295
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 25295
$ws.Range("N68").Value = -26793
$ws.Range("H71").Value = 25295
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 25295
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 75885
$ws.Range("N71").Value = -83373
$ws.Range("H113").Value = 5796
$ws.Range("I113").Value = 6439
$ws.Range("J113").Value = 5428.5713
$ws.Range("K113").Value = 6439
$ws.Range("L113").Value = 5428.5713
$ws.Range("M113").Value = -3185
$ws.Range("N113").Value = -11936.5713
$ws.Range("H116").Value = 5212.7144
$ws.Range("I116").Value = 5999
$ws.Range("J116").Value = 4998.273
$ws.Range("K116").Value = 5999
$ws.Range("L116").Value = 4998.273
$ws.Range("M116").Value = -2557
$ws.Range("N116").Value = -11882.273
$ws.Range("H132").Value = 2564.8914
$ws.Range("I132").Value = 975.53125
$ws.Range("J132").Value = 6197.7144
$ws.Range("K132").Value = 2926.59375
$ws.Range("L132").Value = 18593.1432
$ws.Range("M132").Value = -396.59375
$ws.Range("N132").Value = -23653.1432
$ws.Range("H138").Value = 1816.4517
$ws.Range("I138").Value = 1412.8334
$ws.Range("J138").Value = 2375.3076
$ws.Range("K138").Value = 4238.5002
$ws.Range("L138").Value = 7125.9228
$ws.Range("M138").Value = 901.4997999999996
$ws.Range("N138").Value = -17405.9228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 837552.9399999999
$ws.Range("I45").Value = 1670500.9
$ws.Range("J45").Value = 4605
$ws.Range("K45").Value = 1670500.9
$ws.Range("L45").Value = 4605
$ws.Range("M45").Value = -1670123.9
$ws.Range("N45").Value = -5359

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 19607
$ws.Range("I75").Value = 19607
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 19607
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -18671
$ws.Range("H78").Value = 19607
$ws.Range("I78").Value = 19607
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 58821
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -54141
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").ClearContents()
$ws.Range("N87").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").ClearContents()
$ws.Range("N90").Value = 0
$ws.Range("H113").Value = 5875
$ws.Range("I113").Value = 5875
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 5875
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -3705

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1240.3928
$ws.Range("I31").Value = 1240.3928
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1240.3928
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -945.3928000000001
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 1240.3928
$ws.Range("I34").Value = 1240.3928
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1240.3928
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1038.3928
$ws.Range("N34").ClearContents()
$ws.Range("H62").Value = 36166.668
$ws.Range("I62").Value = 36166.668
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 36166.668
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -35542.668
$ws.Range("H65").Value = 36166.668
$ws.Range("I65").Value = 36166.668
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 180833.34
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -177713.34
$ws.Range("H94").Value = 2681.875
$ws.Range("I94").Value = 3573.6667
$ws.Range("J94").Value = 2146.8
$ws.Range("K94").Value = 3573.6667
$ws.Range("L94").Value = 2146.8
$ws.Range("M94").Value = -3122.6667
$ws.Range("N94").Value = -3048.8
$ws.Range("H134").Value = 2125.889
$ws.Range("I134").Value = 2058.4666
$ws.Range("J134").Value = 2463
$ws.Range("K134").Value = 6175.399800000001
$ws.Range("L134").Value = 7389
$ws.Range("M134").Value = -3640.399800000001
$ws.Range("N134").Value = -12459

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 8559.9
$ws.Range("I9").Value = 862.6
$ws.Range("J9").Value = 16257.2
$ws.Range("K9").Value = 2587.8
$ws.Range("L9").Value = 48771.60000000001
$ws.Range("M9").Value = -2363.8
$ws.Range("N9").Value = -49219.60000000001
$ws.Range("H12").Value = 360.5
$ws.Range("I12").Value = 141
$ws.Range("J12").Value = 470.25
$ws.Range("K12").Value = 423
$ws.Range("L12").Value = 1410.75
$ws.Range("M12").Value = -250
$ws.Range("N12").Value = -1756.75
$ws.Range("H113").Value = 376.75
$ws.Range("I113").Value = 500
$ws.Range("J113").Value = 253.5
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 760.5
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -5100.5
$ws.Range("H129").Value = 113655
$ws.Range("I129").Value = 222600.89
$ws.Range("J129").Value = 4709.1113
$ws.Range("K129").Value = 667802.67
$ws.Range("L129").Value = 14127.3339
$ws.Range("M129").Value = -662802.67
$ws.Range("N129").Value = -24127.3339
$ws.Range("H131").Value = 71041.625
$ws.Range("I131").Value = 158280.86
$ws.Range("J131").Value = 3188.889
$ws.Range("K131").Value = 474842.58
$ws.Range("L131").Value = 9566.667000000001
$ws.Range("M131").Value = -469802.58
$ws.Range("N131").Value = -19646.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 14800
$ws.Range("I20").Value = 14000
$ws.Range("J20").Value = 16000
$ws.Range("K20").Value = 14000
$ws.Range("L20").Value = 16000
$ws.Range("M20").Value = -13755
$ws.Range("N20").Value = -16490
$ws.Range("H70").Value = 5262.8335
$ws.Range("I70").Value = 5325.4
$ws.Range("J70").Value = 4950
$ws.Range("K70").Value = 5325.4
$ws.Range("L70").Value = 4950
$ws.Range("M70").Value = -5055.4
$ws.Range("N70").Value = -5490
$ws.Range("H73").Value = 5262.8335
$ws.Range("I73").Value = 5325.4
$ws.Range("J73").Value = 4950
$ws.Range("K73").Value = 5325.4
$ws.Range("L73").Value = 4950
$ws.Range("M73").Value = -4389.4
$ws.Range("N73").Value = -6822
$ws.Range("H102").Value = 3752.75
$ws.Range("I102").Value = 3506
$ws.Range("J102").Value = 3999.5
$ws.Range("K102").Value = 3506
$ws.Range("L102").Value = 3999.5
$ws.Range("M102").Value = -1884
$ws.Range("N102").Value = -7243.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 2500
$ws.Range("I4").Value = 2500
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2500
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -2387
$ws.Range("H28").Value = 2500
$ws.Range("I28").Value = 2500
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 2500
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -2268
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H37").Value = 2500
$ws.Range("I37").Value = 2500
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 2500
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -2393
$ws.Range("H68").Value = 2575.1428
$ws.Range("I68").Value = 1667
$ws.Range("J68").Value = 3256.25
$ws.Range("K68").Value = 1667
$ws.Range("L68").Value = 3256.25
$ws.Range("M68").Value = -918
$ws.Range("N68").Value = -4754.25
$ws.Range("H71").Value = 2575.1428
$ws.Range("I71").Value = 1667
$ws.Range("J71").Value = 3256.25
$ws.Range("K71").Value = 8335
$ws.Range("L71").Value = 16281.25
$ws.Range("M71").Value = -4591
$ws.Range("N71").Value = -23769.25
$ws.Range("H101").Value = 27908
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 27908
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 27908
$ws.Range("N101").Value = -34398
$ws.Range("H136").Value = 4148.5
$ws.Range("I136").Value = 3498.8333
$ws.Range("J136").Value = 6097.5
$ws.Range("K136").Value = 10496.4999
$ws.Range("L136").Value = 18292.5
$ws.Range("M136").Value = -7946.499899999999
$ws.Range("N136").Value = -23392.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 10999.8
$ws.Range("I29").Value = 11250
$ws.Range("J29").Value = 9999
$ws.Range("K29").Value = 11250
$ws.Range("L29").Value = 9999
$ws.Range("M29").Value = -10960
$ws.Range("N29").Value = -10579
$ws.Range("H62").Value = 6115.75
$ws.Range("I62").Value = 6515.364
$ws.Range("J62").Value = 5236.6
$ws.Range("K62").Value = 6515.364
$ws.Range("L62").Value = 5236.6
$ws.Range("M62").Value = -5891.364
$ws.Range("N62").Value = -6484.6
$ws.Range("H65").Value = 6115.75
$ws.Range("I65").Value = 6515.364
$ws.Range("J65").Value = 5236.6
$ws.Range("K65").Value = 32576.82
$ws.Range("L65").Value = 26183
$ws.Range("M65").Value = -29456.82
$ws.Range("N65").Value = -32423
